$wb = $excel.ActiveWorkbook

$ws_ALC = $wb.Worksheets.Item("ALC")
$ws_ARM = $wb.Worksheets.Item("ARM")
$ws_BSM = $wb.Worksheets.Item("BSM")
$ws_CRP = $wb.Worksheets.Item("CRP")
$ws_CUL = $wb.Worksheets.Item("CUL")
$ws_GSM = $wb.Worksheets.Item("GSM")
$ws_LTW = $wb.Worksheets.Item("LTW")
$ws_WVR = $wb.Worksheets.Item("WVR")

# ALC row 40
$ws_ALC.Range("H40").Value = 2015.174
$ws_ALC.Range("I40").Value = 1857.5883
$ws_ALC.Range("J40").Value = 2461.6667
$ws_ALC.Range("K40").Value = 1857.5883
$ws_ALC.Range("L40").Value = 2461.6667
$ws_ALC.Range("M40").Value = -1682.5883
$ws_ALC.Range("N40").Value = -2811.6667

# ALC row 64
$ws_ALC.Range("H64").Value = 2994.2856
$ws_ALC.Range("I64").Value = 2994.2856
$ws_ALC.Range("J64").Value = 0
$ws_ALC.Range("K64").Value = 2994.2856
$ws_ALC.Range("L64").Value = 0
$ws_ALC.Range("M64").Value = -2746.2856
$ws_ALC.Range("N64").ClearContents()

# ALC row 67
$ws_ALC.Range("H67").Value = 2994.2856
$ws_ALC.Range("I67").Value = 2994.2856
$ws_ALC.Range("J67").Value = 0
$ws_ALC.Range("K67").Value = 2994.2856
$ws_ALC.Range("L67").Value = 0
$ws_ALC.Range("M67").Value = -2136.2856
$ws_ALC.Range("N67").ClearContents()

# ALC row 76
$ws_ALC.Range("H76").Value = 2851.6086
$ws_ALC.Range("I76").Value = 2684.0527
$ws_ALC.Range("J76").Value = 3647.5
$ws_ALC.Range("K76").Value = 2684.0527
$ws_ALC.Range("L76").Value = 3647.5
$ws_ALC.Range("M76").Value = -2369.0527
$ws_ALC.Range("N76").Value = -4277.5

# ALC row 79
$ws_ALC.Range("H79").Value = 2851.6086
$ws_ALC.Range("I79").Value = 2684.0527
$ws_ALC.Range("J79").Value = 3647.5
$ws_ALC.Range("K79").Value = 2684.0527
$ws_ALC.Range("L79").Value = 3647.5
$ws_ALC.Range("M79").Value = -1592.0527
$ws_ALC.Range("N79").Value = -5831.5

# ALC row 132
$ws_ALC.Range("H132").Value = 1833233.1
$ws_ALC.Range("I132").Value = 2748230.8
$ws_ALC.Range("J132").Value = 3238.2307
$ws_ALC.Range("K132").Value = 8244692.399999999
$ws_ALC.Range("L132").Value = 9714.6921
$ws_ALC.Range("M132").Value = -8242162.399999999
$ws_ALC.Range("N132").Value = -14774.6921

# ARM row 9
$ws_ARM.Range("H9").Value = 7050
$ws_ARM.Range("I9").Value = 3000
$ws_ARM.Range("J9").Value = 8400
$ws_ARM.Range("K9").Value = 3000
$ws_ARM.Range("L9").Value = 8400
$ws_ARM.Range("M9").Value = -2830
$ws_ARM.Range("N9").Value = -8740

# ARM row 11
$ws_ARM.Range("H11").Value = 5375
$ws_ARM.Range("I11").Value = 1500
$ws_ARM.Range("J11").Value = 6666.6665
$ws_ARM.Range("K11").Value = 1500
$ws_ARM.Range("L11").Value = 6666.6665
$ws_ARM.Range("M11").Value = -1356
$ws_ARM.Range("N11").Value = -6954.6665

# ARM row 13
$ws_ARM.Range("H13").Value = 4003650.5
$ws_ARM.Range("I13").Value = 20000000
$ws_ARM.Range("J13").Value = 4563
$ws_ARM.Range("K13").Value = 20000000
$ws_ARM.Range("L13").Value = 4563
$ws_ARM.Range("M13").Value = -19999856
$ws_ARM.Range("N13").Value = -4851

# ARM row 17
$ws_ARM.Range("H17").Value = 4800
$ws_ARM.Range("I17").Value = 1000
$ws_ARM.Range("J17").Value = 20000
$ws_ARM.Range("K17").Value = 1000
$ws_ARM.Range("L17").Value = 20000
$ws_ARM.Range("M17").Value = -827
$ws_ARM.Range("N17").Value = -20346

# ARM row 20
$ws_ARM.Range("H20").Value = 7050
$ws_ARM.Range("I20").Value = 3000
$ws_ARM.Range("J20").Value = 8400
$ws_ARM.Range("K20").Value = 3000
$ws_ARM.Range("L20").Value = 8400
$ws_ARM.Range("M20").Value = -2730
$ws_ARM.Range("N20").Value = -8940

# ARM row 33
$ws_ARM.Range("H33").Value = 1506412.6
$ws_ARM.Range("I33").Value = 3001925.2
$ws_ARM.Range("J33").Value = 10900
$ws_ARM.Range("K33").Value = 3001925.2
$ws_ARM.Range("L33").Value = 10900
$ws_ARM.Range("M33").Value = -3001596.2
$ws_ARM.Range("N33").Value = -11558

# ARM row 36
$ws_ARM.Range("H36").Value = 4853.273
$ws_ARM.Range("I36").Value = 983.7143
$ws_ARM.Range("J36").Value = 11625
$ws_ARM.Range("K36").Value = 983.7143
$ws_ARM.Range("L36").Value = 11625
$ws_ARM.Range("M36").Value = -637.7143
$ws_ARM.Range("N36").Value = -12317

# ARM row 63
$ws_ARM.Range("H63").Value = 1944.0625
$ws_ARM.Range("I63").Value = 1709.5454
$ws_ARM.Range("J63").Value = 2460
$ws_ARM.Range("K63").Value = 1709.5454
$ws_ARM.Range("L63").Value = 2460
$ws_ARM.Range("M63").Value = -1023.5454
$ws_ARM.Range("N63").Value = -3832

# ARM row 66
$ws_ARM.Range("H66").Value = 1944.0625
$ws_ARM.Range("I66").Value = 1709.5454
$ws_ARM.Range("J66").Value = 2460
$ws_ARM.Range("K66").Value = 8547.726999999999
$ws_ARM.Range("L66").Value = 12300
$ws_ARM.Range("M66").Value = -5115.726999999999
$ws_ARM.Range("N66").Value = -19164

# BSM row 33
$ws_BSM.Range("H33").Value = 30400
$ws_BSM.Range("I33").Value = 1200
$ws_BSM.Range("J33").Value = 45000
$ws_BSM.Range("K33").Value = 1200
$ws_BSM.Range("L33").Value = 45000
$ws_BSM.Range("M33").Value = -864
$ws_BSM.Range("N33").Value = -45672

# CRP row 17
$ws_CRP.Range("H17").Value = 24857
$ws_CRP.Range("I17").Value = 15499.5
$ws_CRP.Range("J17").Value = 28600
$ws_CRP.Range("K17").Value = 15499.5
$ws_CRP.Range("L17").Value = 28600
$ws_CRP.Range("M17").Value = -15325.5
$ws_CRP.Range("N17").Value = -28948

# CRP row 32
$ws_CRP.Range("H32").Value = 0
$ws_CRP.Range("I32").Value = 0
$ws_CRP.Range("J32").Value = 0
$ws_CRP.Range("K32").Value = 0
$ws_CRP.Range("L32").Value = 0
$ws_CRP.Range("M32").ClearContents()

# CUL row 3
$ws_CUL.Range("H3").Value = 4349.3
$ws_CUL.Range("I3").Value = 3776
$ws_CUL.Range("J3").Value = 4922.6
$ws_CUL.Range("K3").Value = 11328
$ws_CUL.Range("L3").Value = 14767.8
$ws_CUL.Range("M3").Value = -11216
$ws_CUL.Range("N3").Value = -14991.8

# CUL row 86
$ws_CUL.Range("H86").Value = 1875.0588
$ws_CUL.Range("I86").Value = 450
$ws_CUL.Range("J86").Value = 2313.5386
$ws_CUL.Range("K86").Value = 1350
$ws_CUL.Range("L86").Value = 6940.6158
$ws_CUL.Range("M86").Value = -164
$ws_CUL.Range("N86").Value = -9312.6158

# CUL row 89
$ws_CUL.Range("H89").Value = 1875.0588
$ws_CUL.Range("I89").Value = 450
$ws_CUL.Range("J89").Value = 2313.5386
$ws_CUL.Range("K89").Value = 4050
$ws_CUL.Range("L89").Value = 20821.8474
$ws_CUL.Range("M89").Value = 1878
$ws_CUL.Range("N89").Value = -32677.8474

# CUL row 113
$ws_CUL.Range("H113").Value = 3663401.8
$ws_CUL.Range("I113").Value = 382.45834
$ws_CUL.Range("J113").Value = 9524232
$ws_CUL.Range("K113").Value = 1147.37502
$ws_CUL.Range("L113").Value = 28572696
$ws_CUL.Range("M113").Value = 1022.62498
$ws_CUL.Range("N113").Value = -28577036

# CUL row 125
$ws_CUL.Range("H125").Value = 15000
$ws_CUL.Range("I125").Value = 2000
$ws_CUL.Range("J125").Value = 28000
$ws_CUL.Range("K125").Value = 6000
$ws_CUL.Range("L125").Value = 84000
$ws_CUL.Range("M125").Value = -1080
$ws_CUL.Range("N125").Value = -93840

# GSM row 70
$ws_GSM.Range("H70").Value = 4231.625
$ws_GSM.Range("I70").Value = 4425.75
$ws_GSM.Range("J70").Value = 4134.5625
$ws_GSM.Range("K70").Value = 4425.75
$ws_GSM.Range("L70").Value = 4134.5625
$ws_GSM.Range("M70").Value = -4155.75
$ws_GSM.Range("N70").Value = -4674.5625

# GSM row 73
$ws_GSM.Range("H73").Value = 4231.625
$ws_GSM.Range("I73").Value = 4425.75
$ws_GSM.Range("J73").Value = 4134.5625
$ws_GSM.Range("K73").Value = 4425.75
$ws_GSM.Range("L73").Value = 4134.5625
$ws_GSM.Range("M73").Value = -3489.75
$ws_GSM.Range("N73").Value = -6006.5625

# GSM row 80
$ws_GSM.Range("H80").Value = 2543
$ws_GSM.Range("I80").Value = 2262.875
$ws_GSM.Range("J80").Value = 2916.5
$ws_GSM.Range("K80").Value = 2262.875
$ws_GSM.Range("L80").Value = 2916.5
$ws_GSM.Range("M80").Value = -1264.875
$ws_GSM.Range("N80").Value = -4912.5

# GSM row 83
$ws_GSM.Range("H83").Value = 2543
$ws_GSM.Range("I83").Value = 2262.875
$ws_GSM.Range("J83").Value = 2916.5
$ws_GSM.Range("K83").Value = 11314.375
$ws_GSM.Range("L83").Value = 14582.5
$ws_GSM.Range("M83").Value = -6322.375
$ws_GSM.Range("N83").Value = -24566.5

# LTW row 30
$ws_LTW.Range("H30").Value = 25175
$ws_LTW.Range("I30").Value = 700
$ws_LTW.Range("J30").Value = 33333.332
$ws_LTW.Range("K30").Value = 700
$ws_LTW.Range("L30").Value = 33333.332
$ws_LTW.Range("M30").Value = -592
$ws_LTW.Range("N30").Value = -33549.332

# WVR row 10
$ws_WVR.Range("H10").Value = 20000000
$ws_WVR.Range("I10").Value = 20000000
$ws_WVR.Range("J10").Value = 0
$ws_WVR.Range("K10").Value = 20000000
$ws_WVR.Range("L10").Value = 0
$ws_WVR.Range("M10").Value = -19999831

# WVR row 109
$ws_WVR.Range("H109").Value = 32166.666
$ws_WVR.Range("I109").Value = 0
$ws_WVR.Range("J109").Value = 32166.666
$ws_WVR.Range("K109").Value = 0
$ws_WVR.Range("L109").Value = 32166.666
$ws_WVR.Range("N109").Value = -34940.666
